$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look like plain numbers
# (e.g. "581.36", "1.00"). Force text format so Excel doesn't coerce
# them into floating point numbers (losing the original formatting /
# exact digit string).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.685.32"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.421.03"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "581.36"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "129.83"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").Value = "7.60"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "0.384"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "4.002.14"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "3.425.17"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "63.670.06"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "25.07"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "9.87"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "5.67"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "384.74"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "0.564"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").Value = "3.557.64"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "73.70"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0000110"
$ws.Range("E26").Value = "  -4.46%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.990"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.155"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.42"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D33").Value = "3.447.92"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "22.93"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.19"
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "6.79"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "163.73"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0776"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "0.786"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "41.38"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.33"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.61"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "23.46"
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "6.73"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "0.898"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.297.40"
$ws.Range("E50").Value = "  -6.94%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0254"
$ws.Range("E51").Value = "  -2.08%  "
